# Commit: "completato operazione con un quadrante"
# Fills in the "hm mozzo [mm]" (column D) measurements for the Q1 quadrant
# block (rows 29-37) on the "Euramet" sheet, updates a few Laumas/Torsiometro
# readings (columns F/G) that were re-measured for that same block, and clears
# out the still-unmeasured rows (38-46) that had been accidentally pre-filled.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Euramet")

# --- Rows 29-37: enter the "hm mozzo" measurement (837) and correct the
#     Laumas [N] / Torsiometro [Nm] readings that changed with it ---

# Row 29
$ws.Range("D29").Value = 837
$ws.Range("F29").Value = -156.5579745837849

# Row 30
$ws.Range("D30").Value = 837
$ws.Range("F30").Value = -157.0411782090434
$ws.Range("G30").Value = -1.3734

# Row 31
$ws.Range("D31").Value = 837

# Row 32
$ws.Range("D32").Value = 837
$ws.Range("G32").Value = -1.3734

# Row 33
$ws.Range("D33").Value = 837

# Row 34
$ws.Range("D34").Value = 837
$ws.Range("F34").Value = -157.0411782090434

# Row 35
$ws.Range("D35").Value = 837
$ws.Range("F35").Value = -156.5579745837849

# Row 36
$ws.Range("D36").Value = 837
$ws.Range("F36").Value = -156.5579745837849

# Row 37
$ws.Range("D37").Value = 837
$ws.Range("F37").Value = -157.0411782090434
$ws.Range("G37").Value = -1.1772

# --- Rows 38-46: these had been pre-filled with placeholder readings before
#     the quadrant was actually measured; clear E:H back out since only one
#     quadrant's worth of steps (29-37) was completed ---
$ws.Range("E38:H46").ClearContents()
